$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: clear the polite_expressions ("nan") value to an empty string,
# keeping it as text (matches the target inlineStr with no content).
$ws.Range("C10").Value = "'"
$ws.Range("C10").Style = "Normal"

# New row 11: parisk annotation row.
$ws.Range("A11").Value = "parisk"

# politeness_score (B11) must stay text "3" (not numeric 3), so force a
# text format before assigning, then drop the custom format again so no
# stray style is left behind on the cell.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "3"
$ws.Range("B11").Style = "Normal"

$ws.Range("C11").Value = "nan"
$ws.Range("D11").Value = "SUG"
$ws.Range("E11").Value = "MET"
$ws.Range("F11").Value = "d3fb2dcb-ee08-4432-9f4b-c252dbb3433f"
$ws.Range("G11").Value = "SJ3dBGZ0Z_annotated.xlsx"
$ws.Range("H11").Value = "We evaluate our method on NLP task for two reasons: 1) they are particularly well-suited for evaluating our method (naturally large output spaces) 2) we did not dispose of the computational resources to tackle tasks from other domains such as vision (e.g. Flickr100M) which requires hundreds of GPUs for weeks."
